$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -1.692275877396497
$ws.Range("C2").Value = -0.003275877396490101
$ws.Range("D2").Value = 0.3477207439705694
$ws.Range("E2").Value = 0.5877241226035039
$ws.Range("F2").Value = 0.09872412260351376
$ws.Range("G2").Value = -0.3522758773964938
$ws.Range("H2").Value = -0.2722758773964955
$ws.Range("B3").Value = 0.1640906121094658
$ws.Range("C3").Value = 0.5150872334765253
$ws.Range("D3").Value = 0.7550906121094598
$ws.Range("E3").Value = 0.2660906121094697
$ws.Range("F3").Value = -0.1849093878905379
$ws.Range("G3").Value = -0.1049093878905396
$ws.Range("B4").Value = 0.7443053338551252
$ws.Range("C4").Value = 0.9843087124880596
$ws.Range("D4").Value = 0.4953087124880696
$ws.Range("E4").Value = 0.04430871248806201
$ws.Range("F4").Value = 0.1243087124880603
$ws.Range("G4").Value = 0.1443087124880705
$ws.Range("H4").Value = 0.2443087124880648
$ws.Range("I4").Value = 0.1443087124880705
$ws.Range("J4").Value = -0.07569287514151168
$ws.Range("B5").Value = 2.220230824935162
$ws.Range("C5").Value = 1.731230824935172
$ws.Range("D5").Value = 1.280230824935164
$ws.Range("E5").Value = 1.360230824935162
$ws.Range("F5").Value = 1.380230824935173
$ws.Range("G5").Value = 1.480230824935167
$ws.Range("H5").Value = 1.380230824935173
$ws.Range("I5").Value = 1.16022923730559
$ws.Range("B6").Value = 0.3925192495140906
$ws.Range("C6").Value = -0.058480750485917
$ws.Range("D6").Value = 0.02151924951408128
$ws.Range("E6").Value = 0.04151924951409147
$ws.Range("F6").Value = 0.1415192495140858
$ws.Range("G6").Value = 0.04151924951409147
$ws.Range("H6").Value = -0.1784823381154907
$ws.Range("B7").Value = 0.1132553379850393
$ws.Range("C7").Value = 0.1932553379850376
$ws.Range("D7").Value = 0.2132553379850478
$ws.Range("E7").Value = 0.3132553379850421
$ws.Range("F7").Value = 0.2132553379850478
$ws.Range("G7").Value = -0.006746249644534352
$ws.Range("B8").Value = 0.8750144670621453
$ws.Range("C8").Value = 0.8950144670621555
$ws.Range("D8").Value = 0.9950144670621498
$ws.Range("E8").Value = 0.8950144670621555
$ws.Range("F8").Value = 0.6750128794325734
$ws.Range("G8").Value = 0.895016752127448
$ws.Range("H8").Value = 0.7950110251378547
$ws.Range("I8").Value = 0.9550144670621578
$ws.Range("B9").Value = 0.4314358256754829
$ws.Range("C9").Value = 0.5314358256754772
$ws.Range("D9").Value = 0.4314358256754829
$ws.Range("E9").Value = 0.2114342380459007
$ws.Range("F9").Value = 0.4314381107407754
$ws.Range("G9").Value = 0.3314323837511821
$ws.Range("H9").Value = 0.4914358256754852
$ws.Range("B10").Value = -0.02922200567268469
$ws.Range("C10").Value = -0.129222005672679
$ws.Range("D10").Value = -0.3492235933022612
$ws.Range("E10").Value = -0.1292197206073865
$ws.Range("F10").Value = -0.2292254475969798
$ws.Range("G10").Value = -0.06922200567267674
$ws.Range("B11").Value = -0.02670664067814384
$ws.Range("C11").Value = -0.246708228307726
$ws.Range("D11").Value = -0.0267043556128513
$ws.Range("E11").Value = -0.1267100826024446
$ws.Range("F11").Value = 0.03329335932185844
$ws.Range("B12").Value = -0.3283108325625718
$ws.Range("C12").Value = -0.1083069598676971
$ws.Range("D12").Value = -0.2083126868572904
$ws.Range("E12").Value = -0.04830924493298736
$ws.Range("B13").Value = -0.1236027424349175
$ws.Range("C13").Value = -0.2236084694245108
$ws.Range("D13").Value = -0.06360502750020772
$ws.Range("B14").Value = -0.3154970214097745
$ws.Range("C14").Value = -0.1554935794854714
$ws.Range("B15").Value = 0.08316376585921856
